$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the raw measurement data in B2:B31 with the new recorded values.
# ---------------------------------------------------------------------------
$newVals = @(
    168.11712,
    185.24016,
    182.12688,
    179.0136,
    183.68351999999999,
    171.2304,
    180.57023999999899,
    180.57023999999899,
    180.57023999999899,
    179.0136,
    180.57023999999899,
    179.0136,
    177.45696000000001,
    177.45696000000001,
    179.0136,
    177.45696000000001,
    179.0136,
    179.0136,
    182.12688,
    177.45696000000001,
    180.57023999999899,
    171.2304,
    179.0136,
    180.57023999999899,
    180.57023999999899,
    180.57023999999899,
    180.57023999999899,
    182.12688,
    179.0136,
    179.0136
)

for ($i = 0; $i -lt $newVals.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $newVals[$i]
}

# ---------------------------------------------------------------------------
# 2. Add the two new summary columns: "Mean increase" and "Median increase".
# ---------------------------------------------------------------------------
$ws.Range("D18").Value = "Mean increase"
$ws.Range("D18").Font.Bold = $true

$ws.Range("F18").Value = "Median increase"
$ws.Range("F18").Font.Bold = $true

$ws.Range("D19").Formula = "=((E3 / 114.202998) * 100) - 100"
$ws.Range("D19").Style = "Standaard"

$ws.Range("F19").Formula = "=((E10 / 113.658804) * 100) - 100"
$ws.Range("F19").Style = "Standaard"

# ---------------------------------------------------------------------------
# 3. Update the sheet selection / active cell.
# ---------------------------------------------------------------------------
[void]$ws.Range("F23").Select()

# ---------------------------------------------------------------------------
# 4. Reposition the workbook window (bookViews / workbookView).
# ---------------------------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.Left = 60
$win.Top = 400
